# Generate Report for Handoff
#
# A new source file (c62b932a-7f19-4e07-a4a2-9772090359a3.md) has finished
# handoff, so it gets appended as a third data row to each of the three
# report tables: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$repoBlob = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd1ff4cc51ac351fe7251018796d3ce5fe1e4bd6/e2e/"
$dateFmt  = "yyyy-mm-dd HH:mm:ss"

# Helper: write a value that must stay TEXT (e.g. "True"/"False") even
# though it looks like another type to Excel's auto-detection, then drop
# back to the plain "Normal" style so no stray quote-prefix formatting
# lingers on the cell.
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.ListRows.Add() | Out-Null

$rOverview = $tblOverview.Range.Row + $tblOverview.ListRows.Count

$wsOverview.Cells.Item($rOverview, 1).Value = "c62b932a-7f19-4e07-a4a2-9772090359a3.md"
$wsOverview.Cells.Item($rOverview, 2).Value = "e2e\c62b932a-7f19-4e07-a4a2-9772090359a3.md"
$wsOverview.Cells.Item($rOverview, 3).Value = ".md"
$wsOverview.Cells.Item($rOverview, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item($rOverview, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item($rOverview, 7).Value = "2016-09-05 04:44:48"
$wsOverview.Cells.Item($rOverview, 7).NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item($rOverview, 2),
    ($repoBlob + "c62b932a-7f19-4e07-a4a2-9772090359a3.md"),
    "",
    "",
    "e2e\c62b932a-7f19-4e07-a4a2-9772090359a3.md"
) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$tblZhCn = $wsZhCn.ListObjects.Item(1)
$tblZhCn.ListRows.Add() | Out-Null

$rZhCn = $tblZhCn.Range.Row + $tblZhCn.ListRows.Count

$wsZhCn.Cells.Item($rZhCn, 1).Value = "c62b932a-7f19-4e07-a4a2-9772090359a3.md"
$wsZhCn.Cells.Item($rZhCn, 2).Value = ".md"
$wsZhCn.Cells.Item($rZhCn, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item($rZhCn, 4).Value = "e2e"
$wsZhCn.Cells.Item($rZhCn, 5).Value = "ht"
Set-TextValue $wsZhCn.Cells.Item($rZhCn, 6) "False"
$wsZhCn.Cells.Item($rZhCn, 7).Value = "c62b932a-7f19-4e07-a4a2-9772090359a3.76d34245c170fff639a9763afd00dd7c1b37bac8.zh-cn.xlf"
$wsZhCn.Cells.Item($rZhCn, 8).Value = "2016-09-05 04:44:43"
$wsZhCn.Cells.Item($rZhCn, 8).NumberFormat = $dateFmt
$wsZhCn.Cells.Item($rZhCn, 11).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item($rZhCn, 11).NumberFormat = $dateFmt
Set-TextValue $wsZhCn.Cells.Item($rZhCn, 13) "True"
Set-TextValue $wsZhCn.Cells.Item($rZhCn, 15) "False"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item($rZhCn, 1),
    ($repoBlob + "c62b932a-7f19-4e07-a4a2-9772090359a3.md"),
    "",
    "",
    "c62b932a-7f19-4e07-a4a2-9772090359a3.md"
) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$tblDeDe = $wsDeDe.ListObjects.Item(1)
$tblDeDe.ListRows.Add() | Out-Null

$rDeDe = $tblDeDe.Range.Row + $tblDeDe.ListRows.Count

$wsDeDe.Cells.Item($rDeDe, 1).Value = "c62b932a-7f19-4e07-a4a2-9772090359a3.md"
$wsDeDe.Cells.Item($rDeDe, 2).Value = ".md"
$wsDeDe.Cells.Item($rDeDe, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item($rDeDe, 4).Value = "e2e"
$wsDeDe.Cells.Item($rDeDe, 5).Value = "ht"
Set-TextValue $wsDeDe.Cells.Item($rDeDe, 6) "False"
$wsDeDe.Cells.Item($rDeDe, 7).Value = "c62b932a-7f19-4e07-a4a2-9772090359a3.76d34245c170fff639a9763afd00dd7c1b37bac8.de-de.xlf"
$wsDeDe.Cells.Item($rDeDe, 8).Value = "2016-09-05 04:44:48"
$wsDeDe.Cells.Item($rDeDe, 8).NumberFormat = $dateFmt
$wsDeDe.Cells.Item($rDeDe, 11).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item($rDeDe, 11).NumberFormat = $dateFmt
Set-TextValue $wsDeDe.Cells.Item($rDeDe, 13) "True"
Set-TextValue $wsDeDe.Cells.Item($rDeDe, 15) "False"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item($rDeDe, 1),
    ($repoBlob + "c62b932a-7f19-4e07-a4a2-9772090359a3.md"),
    "",
    "",
    "c62b932a-7f19-4e07-a4a2-9772090359a3.md"
) | Out-Null
